$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Set row heights (w:trHeight val in twips = Height in points * 20)
# Row 1 (config)      -> 517 twips = 25.85 pt
# Row 2 (help)        -> 505 twips = 25.25 pt
# Row 3 (Init)        -> 517 twips = 25.85 pt
# Row 4 (Status)      -> 505 twips = 25.25 pt
# Row 5 (Add <file>)  -> 517 twips = 25.85 pt
# Row 6 (Commit)      -> 517 twips = 25.85 pt
# Row 7 (Reset)       -> 505 twips = 25.25 pt
# Row 8 (Restore)     -> 517 twips = 25.85 pt
# Row 9 (Clean)       -> 505 twips = 25.25 pt
# Row 10 (Branch)     -> 517 twips = 25.85 pt
# Row 11 (empty)      -> 252 twips = 12.6 pt
$heightsPt = @(25.85, 25.25, 25.85, 25.25, 25.85, 25.85, 25.25, 25.85, 25.25, 25.85, 12.6)

for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $t.Rows.Item($i).Height = $heightsPt[$i - 1]
}

# Fix the "Git add . | git add <arquivo>" cell: merge the "add" and " ." runs
# into a single "add ." run (also drops the now-redundant spellcheck marker
# between them), matching the commit's Git Bash command cleanup.
$cell = $t.Rows.Item(5).Cells.Item(2)
$rng = $cell.Range
$rng.Find.Execute("add .", $true, $false, $false, $false, $false, $true, 1, $false, "add .", 2)
